# Commit: "Add checking of missing data when computing points to reacquire"
#
# The "positioning" worksheet had pre-computed reacquisition coordinates
# (Easting -> column S, Elevation -> column T) duplicated from columns D/E
# for every survey point (rows 13-132). The new "missing data" check means
# those values are no longer blindly copied over - the cells are left
# blank (still formatted, but with no content) wherever the check applies,
# which for this sheet's data is every row. Clear the stored values while
# keeping the existing cell formatting intact.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Positioning")

$ws.Range("S13:T132").ClearContents()
